# Daily attendance processing - 2026-01-11 23:57:13
# Normalize the "Recorded By" (column G) values so that when the "System"
# user is listed alongside other recorders, "System" is always moved to the
# front of the comma-separated list (e.g. "dnasr281@gmail.com, System" ->
# "System, dnasr281@gmail.com").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 7).End(-4162).Row  # xlUp = -4162, column 7 = G

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2

    if ($null -eq $val) { continue }
    if (-not ($val -is [string])) { continue }
    if ($val -notmatch ',') { continue }

    $parts = $val -split ', '

    $sysIndex = -1
    for ($i = 0; $i -lt $parts.Length; $i++) {
        if ($parts[$i].Equals('System')) {
            $sysIndex = $i
            break
        }
    }

    if ($sysIndex -gt 0) {
        $newParts = @('System')
        for ($i = 0; $i -lt $parts.Length; $i++) {
            if ($i -ne $sysIndex) {
                $newParts += $parts[$i]
            }
        }
        $cell.Value = [string]::Join(', ', $newParts)
    }
}
